$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 145 (existing rows 145:168 shift down to 146:169,
# preserving their formatting/styles).
$ws.Rows("145:145").Insert()

# Populate the newly inserted row 145 with the new weekly price record.
$ws.Cells.Item(145, 1).Value  = 7
$ws.Cells.Item(145, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(145, 3).Value  = "Ñuble"
$ws.Cells.Item(145, 4).Value  = 44491
$ws.Cells.Item(145, 5).Value  = 16
$ws.Cells.Item(145, 6).Value  = "Fruta"
$ws.Cells.Item(145, 7).Value  = 100108
$ws.Cells.Item(145, 8).Value  = "Tropicales y subtropicales"
$ws.Cells.Item(145, 9).Value  = 100108005
$ws.Cells.Item(145, 10).Value = "Piña"
$ws.Cells.Item(145, 11).Value = "Caramelo"
$ws.Cells.Item(145, 12).Value = "Segunda"
$ws.Cells.Item(145, 13).Value = 120
$ws.Cells.Item(145, 14).Value = 19500
$ws.Cells.Item(145, 15).Value = 20000
$ws.Cells.Item(145, 16).Value = 19750
$ws.Cells.Item(145, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(145, 18).Value = "Ecuador"
$ws.Cells.Item(145, 19).Value = 1411
$ws.Cells.Item(145, 20).Value = 14
